# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 56
    3  = 3360
    5  = 2422
    8  = 1373
    9  = 1086
    10 = 294
    11 = 508
    14 = 98
    16 = 8472
    18 = 2480
    19 = 250
    23 = 580
    27 = 1994
    30 = 1731
    34 = 25
    35 = 38
    39 = 301
    42 = 405
    43 = 116
    45 = 253
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates4 = @{
    2  = 56
    3  = 3360
    5  = 2422
    8  = 1373
    10 = 1086
    11 = 294
    12 = 508
    14 = 98
    16 = 8472
    18 = 2480
    20 = 250
    24 = 580
    28 = 1994
    30 = 1731
    34 = 25
    35 = 38
    39 = 301
    42 = 405
    47 = 116
    49 = 253
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
